$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header row labels to reflect plans without sticky-ids / levels support
$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Move the active selection to F1 as recorded in the saved view
$ws.Range("F1").Select()
